$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New empty cell D2, sharing B2's date-style formatting
$ws.Range("B2").Copy()
$ws.Range("D2").PasteSpecial(-4122)

# Clear the stray hard-coded "22" total under D34
$ws.Range("D34").ClearContents()

# Row 38 summary cells used to hold hard-coded numbers (165 / 176 / 11).
# Clear them and give D38/F38 a [h]:mm duration format, H38 a text format.
$ws.Range("D38").ClearContents()
$ws.Range("F38").ClearContents()
$ws.Range("H38").ClearContents()
$ws.Range("D38").NumberFormat = "[h]:mm"
$ws.Range("F38").NumberFormat = "[h]:mm"
$ws.Range("H38").NumberFormat = "@"

# Move the active selection
$ws.Range("H32").Select()
